$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column A display text (email addresses) for rows 2-7 ---
# Row 8 keeps its original "jakay34@gmail.com" text (unchanged).
$ws.Range("A2").Value = "jakay11@gmail.com"
$ws.Range("A3").Value = "jakay12@gmail.com"
$ws.Range("A4").Value = "jakay13@gmail.com"
$ws.Range("A5").Value = "jakay14@gmail.com"
$ws.Range("A6").Value = "jakay15@gmail.com"
$ws.Range("A7").Value = "jakay16@gmail.com"

# --- Update column B (password) for rows 2-7 to the numeric password value ---
$ws.Range("B2").Value = 12345678
$ws.Range("B3").Value = 12345678
$ws.Range("B4").Value = 12345678
$ws.Range("B5").Value = 12345678
$ws.Range("B6").Value = 12345678
$ws.Range("B7").Value = 12345678

# --- Fix up the hyperlink targets that changed ---
# Row A5's hyperlink now points at jakay34@gmail.com, row A8's hyperlink now
# points at automation@gmail.com (the two effectively swapped targets).
$ws.Hyperlinks.Item(4).Address = "mailto:jakay34@gmail.com"
$ws.Hyperlinks.Item(7).Address = "mailto:automation@gmail.com"

# --- Update the selected cell shown in the sheet view ---
$ws.Range("D8").Select()

Write-Host "edit complete"
